# SettleGame.xlsx edit script
# Applies the "Game design typo/ordering fixes" commit:
#  - Moves the "Boar Ranch" row (row 42) up to row 38 (ahead of Hunting Camp,
#    Cane Fields, Smoke House, Peach Orchard, which all shift down by one row),
#    and fixes its Cost from "2M, 1F" to "1M, 1F".
#  - Fixes Peach Orchard's effect text (now row 42) to drop the "or +1 Active" option.
#  - Fixes Trade Harbor's effect text (row 31).
#  - Fixes Silver Mine's effect text (row 47).
#  - Updates the selected cell shown in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basics")

# --- Move the Boar Ranch row (42) to before Hunting Camp (38) ---
$ws.Rows(42).Copy()
$ws.Rows(38).Insert()
$ws.Rows(43).Delete()

# Re-apply the normal (non-hidden, "has alt name") row formatting to the moved
# row, since Insert() picks up formatting from the row above it.
$ws.Range("A35:G35").Copy()
$ws.Range("A38:G38").PasteSpecial(-4122)

# Content fix: Boar Ranch cost is 1 Material + 1 Food, not 2 Material + 1 Food.
$ws.Range("C38").Value = "1M, 1F"

# --- Formula/text fixes ---
$ws.Range("G31").Formula = '="-1 Food, -1 Material, +1 Treasure"'
$ws.Range("G42").Formula = '="+2 Food"'
$ws.Range("G47").Formula = '="-2 Food, +1 Treasure"'

# --- Sheet view: update selected cell ---
$ws.Activate()
$ws.Range("G42").Select()
